$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.981.57"
$ws.Range("E2").Value = "  +1.79%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.454.87"
$ws.Range("E3").Value = "  +1.26%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "575.78"
$ws.Range("E5").Value = "  +1.14%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "160.19"
$ws.Range("E6").Value = "  +2.71%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.608"
$ws.Range("E7").Value = "  +5.05%  "
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.460.61"
$ws.Range("E9").Value = "  +1.33%  "
$ws.Range("E10").Value = "  -1.04%  "
$ws.Range("E11").Value = "  +1.33%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.449"
$ws.Range("E12").Value = "  +1.31%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.050.70"
$ws.Range("E13").Value = "  +1.32%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.134"
$ws.Range("E14").Value = "  +0.47%  "
$ws.Range("E15").Value = "  -0.11%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "28.30"
$ws.Range("E16").Value = "  +2.33%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "64.985.75"
$ws.Range("E17").Value = "  +1.96%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.445.77"
$ws.Range("E18").Value = "  +1.13%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.46"
$ws.Range("E19").Value = "  +1.03%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.29"
$ws.Range("E20").Value = "  +0.30%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "381.39"
$ws.Range("E21").Value = "  -0.69%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "8.15"
$ws.Range("E22").Value = "  +0.42%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.555"
$ws.Range("E23").Value = "  +2.99%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "72.96"
$ws.Range("E24").Value = "  -0.41%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.998"
$ws.Range("E25").Value = "  -0.20%  "
$ws.Range("E26").Value = "  +0.59%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.06"
$ws.Range("E27").Value = "  +5.37%  "
$ws.Range("E28").Value = "  -0.29%  "
$ws.Range("E29").Value = "  +0.04%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.52"
$ws.Range("E30").Value = "  +8.28%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.17"
$ws.Range("E31").Value = "  +0.29%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.04"
$ws.Range("E32").Value = "  +2.07%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "23.62"
$ws.Range("E33").Value = "  +0.77%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.28"
$ws.Range("E34").Value = "  +5.02%  "
$ws.Range("E35").Value = "  +11.00%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "161.28"
$ws.Range("E36").Value = "  +0.29%  "
$ws.Range("E37").Value = "  +4.49%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0784"
$ws.Range("E38").Value = "  +2.52%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.88"
$ws.Range("E39").Value = "  +6.95%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "27.07"
$ws.Range("E40").Value = "  +0.48%  "
$ws.Range("B41").Value = "Maker"
$ws.Range("C41").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.903.08"
$ws.Range("E41").Value = "  -0.07%  "
$ws.Range("B42").Value = "Filecoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.71"
$ws.Range("E42").Value = "  +7.49%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "43.13"
$ws.Range("E43").Value = "  +2.25%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0318"
$ws.Range("E44").Value = "  -0.16%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.779"
$ws.Range("E45").Value = "  +2.91%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "25.91"
$ws.Range("E46").Value = "  +11.23%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "322.83"
$ws.Range("E47").Value = "  +10.58%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.09"
$ws.Range("E48").Value = "  +2.47%  "
$ws.Range("E49").Value = "  +1.67%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.876"
$ws.Range("E50").Value = "  +3.98%  "
$ws.Range("E51").Value = "  +0.43%  "
